# Update crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.091.97'
$ws.Range("E2").Value = '  +3.20%  '
$ws.Range("D3").Value = '2.296.70'
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'310.28"
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").Value = "'100.90"
$ws.Range("E6").Value = '  +7.38%  '
$ws.Range("E7").Value = '  +2.71%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = '  +7.66%  '
$ws.Range("D10").Value = "'35.89"
$ws.Range("E10").Value = '  +3.88%  '
$ws.Range("D11").Value = "'0.0823"
$ws.Range("E11").Value = '  +4.80%  '
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("E13").Value = '  +7.82%  '
$ws.Range("D14").Value = '2.654.07'
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").Value = "'14.98"
$ws.Range("E15").Value = '  +5.00%  '
$ws.Range("D16").Value = '2.303.60'
$ws.Range("E16").Value = '  +2.04%  '
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").Value = '43.007.60'
$ws.Range("E18").Value = '  +3.26%  '
$ws.Range("D19").Value = "'12.53"
$ws.Range("E19").Value = '  +2.29%  '
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("D21").Value = "'6.08"
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("D22").Value = "'68.46"
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").Value = "'240.15"
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = '  +5.11%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = "'2.62"
$ws.Range("E25").Value = '  +2.73%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = "'24.46"
$ws.Range("E27").Value = '  +4.01%  '
$ws.Range("D28").Value = "'38.42"
$ws.Range("E28").Value = '  +6.86%  '
$ws.Range("D29").Value = "'2.32"
$ws.Range("E29").Value = '  +10.57%  '
$ws.Range("E30").Value = '  +2.58%  '
$ws.Range("D31").Value = "'167.33"
$ws.Range("E31").Value = '  +5.09%  '
$ws.Range("D32").Value = "'5.32"
$ws.Range("E32").Value = '  +2.70%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("D35").Value = "'17.75"
$ws.Range("E35").Value = '  +5.44%  '
$ws.Range("D36").Value = "'0.0740"
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("E37").Value = '  +3.23%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = '  +2.31%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = "'1.83"
$ws.Range("E40").Value = '  +1.33%  '
$ws.Range("D41").Value = "'4.23"
$ws.Range("E41").Value = '  +7.05%  '
$ws.Range("D42").Value = "'2.30"
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.974.84'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = "'0.0289"
$ws.Range("E44").Value = '  +3.09%  '
$ws.Range("D45").Value = "'19.04"
$ws.Range("E45").Value = '  +2.69%  '
$ws.Range("D46").Value = "'3.03"
$ws.Range("E46").Value = '  +4.31%  '
$ws.Range("D47").Value = "'9.85"
$ws.Range("E47").Value = '  +0.51%  '
$ws.Range("D48").Value = "'55.86"
$ws.Range("E48").Value = '  +6.38%  '
$ws.Range("D49").Value = "'2.94"
$ws.Range("E49").Value = '  +17.35%  '
$ws.Range("D50").Value = '2.523.71'
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("D51").Value = "'1.53"
$ws.Range("E51").Value = '  +2.87%  '
